{"js": "// Resume content overhaul: rename the candidate, swap the contact line,\n// rewrite the summary / skills / certifications / achievements / work\n// history / education / interests / extracurricular / associations\n// sections to match the new \"electrical engineering student\" resume,\n// dropping several bullet points and whole job entries along the way.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Paragraphs whose text changes but whose paragraph (and run formatting)\n// stays in place.\nconst replacements = {\n  0: \"Harrel Ellis\",\n  1: \"4843732 | harel.ellis@gmail.com | Number 4, Bryce Terrace, St Augustine\",\n  4: \"A 22-year-old electrical engineering student with a passion for the field. Seeking opportunities to apply skills in MATLAB, C++ Programming, and Electrical Wiring. Committed to enhancing communication skills and delivering innovative solutions in the electrical engineering industry.\",\n  6: \"\u2022 MATLAB\",\n  7: \"\u2022 C++ Programming\",\n  8: \"\u2022 Electrical Wiring\",\n  9: \"\u2022 Communication Skills\",\n  12: \"\u2022 Electrical and Computer Engineering degree\",\n  17: \"\u2022 CXE examinations\",\n  21: \"Apprentice Plumber\",\n  22: \"Position: Apprentice Plumber\",\n  23: \"Responsible for assisting with plumbing tasks and installations. Successfully completed an electrical wiring project for residential use and a BJT data mining project for programming.\",\n  30: \"Queen's Royal College 2015-2022\",\n  31: \"UWI Engineering 2023-2026\",\n  33: \"\u2022 Electronics\",\n  34: \"\u2022 Robotics\",\n  35: \"\u2022 Programming\",\n  36: \"\u2022 Energy Systems\",\n  39: \"\u2022 Robotics Club: Active member participating in various projects and competitions.\",\n  40: \"\u2022 Football: Engaged in team sports and building teamwork skills.\",\n  45: \"\u2022 IEEE\",\n};\n\n// Paragraphs that are removed entirely (whole bullet points / job entries /\n// the Volunteer Experience section).\nconst deletions = [\n  10, // \"\u2022 Graphic design\"\n  13, // \"\u2022 Python for Data Science\"\n  14, // \"\u2022 Google Analytics Certified\"\n  15, // \"\u2022 ITIL Foundation\"\n  18, // \"\u2022 Developed key automation system\"\n  19, // \"\u2022 Led a successful product launch\"\n  24, // \"BCA Corporation\"\n  25, // \"As a Project helper, ...\"\n  26, // \"DAF Computing\"\n  27, // \"Leader\"\n  28, // \"As a Project leader ensuring ...\"\n  37, // \"\u2022 Photography\"\n  41, // \"Volunteer Experience\"\n  42, // \"\u2022 Habitat for Humanity volunteer, ...\"\n  43, // \"\u2022 helped at the aids foundation ...\"\n  46, // \"\u2022 Association for Computing Machinery\"\n];\n\nfor (const [idx, text] of Object.entries(replacements)) {\n  items[Number(idx)].insertText(text, \"Replace\");\n}\n\nfor (const idx of deletions) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Resume content overhaul: rename the candidate, swap the contact line,\n# rewrite the summary / skills / certifications / achievements / work\n# history / education / interests / extracurricular / associations\n# sections to match the new \"electrical engineering student\" resume,\n# dropping several bullet points and whole job entries along the way.\n\n$d = $word.ActiveDocument\n\n# Paragraphs whose text changes but whose paragraph (and run formatting)\n# stays in place. Keys are 1-based Paragraphs() indices (matching the\n# document's current/original layout \u2014 replacing text never shifts\n# paragraph indices, so these can all be applied in any order).\n$replacements = @{\n    1  = \"Harrel Ellis\"\n    2  = \"4843732 | harel.ellis@gmail.com | Number 4, Bryce Terrace, St Augustine\"\n    5  = \"A 22-year-old electrical engineering student with a passion for the field. Seeking opportunities to apply skills in MATLAB, C++ Programming, and Electrical Wiring. Committed to enhancing communication skills and delivering innovative solutions in the electrical engineering industry.\"\n    7  = \"\u2022 MATLAB\"\n    8  = \"\u2022 C++ Programming\"\n    9  = \"\u2022 Electrical Wiring\"\n    10 = \"\u2022 Communication Skills\"\n    13 = \"\u2022 Electrical and Computer Engineering degree\"\n    18 = \"\u2022 CXE examinations\"\n    22 = \"Apprentice Plumber\"\n    23 = \"Position: Apprentice Plumber\"\n    24 = \"Responsible for assisting with plumbing tasks and installations. Successfully completed an electrical wiring project for residential use and a BJT data mining project for programming.\"\n    31 = \"Queen's Royal College 2015-2022\"\n    32 = \"UWI Engineering 2023-2026\"\n    34 = \"\u2022 Electronics\"\n    35 = \"\u2022 Robotics\"\n    36 = \"\u2022 Programming\"\n    37 = \"\u2022 Energy Systems\"\n    40 = \"\u2022 Robotics Club: Active member participating in various projects and competitions.\"\n    41 = \"\u2022 Football: Engaged in team sports and building teamwork skills.\"\n    46 = \"\u2022 IEEE\"\n}\n\nforeach ($idx in $replacements.Keys) {\n    $d.Paragraphs($idx).Range.Text = $replacements[$idx]\n}\n\n# Paragraphs that are removed entirely (whole bullet points / job entries /\n# the Volunteer Experience section). Delete from the highest index down so\n# earlier deletions never shift the index of a paragraph still pending\n# deletion.\n$deletions = @(11, 14, 15, 16, 19, 20, 25, 26, 27, 28, 29, 38, 42, 43, 44, 47)\n\nforeach ($idx in ($deletions | Sort-Object -Descending)) {\n    $d.Paragraphs($idx).Range.Delete()\n}\n"}
